$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.8268497330282227
$ws.Range("D2").Value = 0.8177987028702338
$ws.Range("E2").Value = 0.8268497330282227
$ws.Range("F2").Value = 0.8123177461656231
$ws.Range("C3").Value = 0.9077040427154843
$ws.Range("D3").Value = 0.903320740687993
$ws.Range("E3").Value = 0.9077040427154843
$ws.Range("F3").Value = 0.9008123345596354
$ws.Range("C4").Value = 0.9944698703279939
$ws.Range("D4").Value = 0.9945004526621831
$ws.Range("E4").Value = 0.9944698703279939
$ws.Range("F4").Value = 0.9917124722746894
$ws.Range("C5").Value = 0.7864225781845919
$ws.Range("D5").Value = 0.7872906718385223
$ws.Range("E5").Value = 0.7864225781845919
$ws.Range("F5").Value = 0.7831637240824391
$ws.Range("C6").Value = 0.9189549961861174
$ws.Range("D6").Value = 0.8992381992068549
$ws.Range("E6").Value = 0.9189549961861174
$ws.Range("F6").Value = 0.8945725637734619
$ws.Range("C7").Value = 0.9551868802440885
$ws.Range("D7").Value = 0.9472949768923256
$ws.Range("E7").Value = 0.9551868802440885
$ws.Range("F7").Value = 0.9419755365761606
$ws.Range("C8").Value = 0.9736842105263158
$ws.Range("D8").Value = 0.964893852277753
$ws.Range("E8").Value = 0.9736842105263158
$ws.Range("F8").Value = 0.9621167640204594
$ws.Range("C9").Value = 0.9794050343249427
$ws.Range("D9").Value = 0.9798291869360995
$ws.Range("E9").Value = 0.9794050343249427
$ws.Range("F9").Value = 0.9692146929273422
$ws.Range("C10").Value = 0.9662471395881007
$ws.Range("D10").Value = 0.9586968623059617
$ws.Range("E10").Value = 0.9662471395881007
$ws.Range("F10").Value = 0.9544755771023257
$ws.Range("C11").Value = 0.9610983981693364
$ws.Range("D11").Value = 0.9576300588995028
$ws.Range("E11").Value = 0.9610983981693364
$ws.Range("F11").Value = 0.9569191488828531
$ws.Range("C12").Value = 0.9464149504195271
$ws.Range("D12").Value = 0.9438230591104911
$ws.Range("E12").Value = 0.9464149504195271
$ws.Range("F12").Value = 0.942051978867804
$ws.Range("C13").Value = 0.9427917620137299
$ws.Range("D13").Value = 0.9386830439089836
$ws.Range("E13").Value = 0.9427917620137299
$ws.Range("F13").Value = 0.934549663266837
$ws.Range("C14").Value = 0.9856979405034325
$ws.Range("D14").Value = 0.9824567065882263
$ws.Range("E14").Value = 0.9856979405034325
$ws.Range("F14").Value = 0.981738016480412
$ws.Range("C15").Value = 0.9792143401983219
$ws.Range("D15").Value = 0.9731808319745648
$ws.Range("E15").Value = 0.9792143401983219
$ws.Range("F15").Value = 0.970918224366996
$ws.Range("C16").Value = 0.9893211289092296
$ws.Range("D16").Value = 0.9894351671970029
$ws.Range("E16").Value = 0.9893211289092296
$ws.Range("F16").Value = 0.9840103559779684
$ws.Range("C17").Value = 0.9679633867276888
$ws.Range("D17").Value = 0.9603865694003006
$ws.Range("E17").Value = 0.9679633867276888
$ws.Range("F17").Value = 0.9555857458440002
$ws.Range("C18").Value = 0.9607170099160945
$ws.Range("D18").Value = 0.9589298282463142
$ws.Range("E18").Value = 0.9607170099160945
$ws.Range("F18").Value = 0.9496149427735556
$ws.Range("C19").Value = 0.8728070175438597
$ws.Range("D19").Value = 0.8409987613237453
$ws.Range("E19").Value = 0.8728070175438597
$ws.Range("F19").Value = 0.8327191554311399
$ws.Range("C20").Value = 0.9364988558352403
$ws.Range("D20").Value = 0.9256431109179933
$ws.Range("E20").Value = 0.9364988558352403
$ws.Range("F20").Value = 0.9100347432995428
$ws.Range("C21").Value = 0.9586193745232647
$ws.Range("D21").Value = 0.9528569127537951
$ws.Range("E21").Value = 0.9586193745232647
$ws.Range("F21").Value = 0.9438677719481351
$ws.Range("C22").Value = 0.9582379862700229
$ws.Range("D22").Value = 0.9539325831904284
$ws.Range("E22").Value = 0.9582379862700229
$ws.Range("F22").Value = 0.9469725454799259
$ws.Range("C23").Value = 0.9811212814645309
$ws.Range("D23").Value = 0.9754349468428717
$ws.Range("E23").Value = 0.9811212814645309
$ws.Range("F23").Value = 0.9746755561320826
$ws.Range("C24").Value = 0.994279176201373
$ws.Range("D24").Value = 0.9943119040263079
$ws.Range("E24").Value = 0.994279176201373
$ws.Range("F24").Value = 0.9914269697291946
$ws.Range("C25").Value = 0.9881769641495042
$ws.Range("D25").Value = 0.9883167483262262
$ws.Range("E25").Value = 0.9881769641495042
$ws.Range("F25").Value = 0.9823006000810917
$ws.Range("C26").Value = 0.9956140350877193
$ws.Range("D26").Value = 0.9956332717759311
$ws.Range("E26").Value = 0.9956140350877193
$ws.Range("F26").Value = 0.9934258723732408
$ws.Range("C27").Value = 0.9872234935163997
$ws.Range("D27").Value = 0.9873867326343252
$ws.Range("E27").Value = 0.9872234935163997
$ws.Range("F27").Value = 0.9808763124334327
$ws.Range("C28").Value = 0.9569031273836766
$ws.Range("D28").Value = 0.9435498674864152
$ws.Range("E28").Value = 0.9569031273836766
$ws.Range("F28").Value = 0.9377883522955645
$ws.Range("C29").Value = 0.8838672768878718
$ws.Range("D29").Value = 0.8831391013379382
$ws.Range("E29").Value = 0.8838672768878718
$ws.Range("F29").Value = 0.8785273189346098
$ws.Range("C30").Value = 0.9509916094584286
$ws.Range("D30").Value = 0.9496097270464978
$ws.Range("E30").Value = 0.9509916094584286
$ws.Range("F30").Value = 0.9423296952066764
$ws.Range("C31").Value = 0.9431731502669718
$ws.Range("D31").Value = 0.9376761152443363
$ws.Range("E31").Value = 0.9431731502669718
$ws.Range("F31").Value = 0.9364232454518965
$ws.Range("C32").Value = 0.9908466819221968
$ws.Range("D32").Value = 0.988001369824934
$ws.Range("E32").Value = 0.9908466819221968
$ws.Range("F32").Value = 0.9868182720113153
$ws.Range("C33").Value = 0.9603356216628528
$ws.Range("D33").Value = 0.9586697165006033
$ws.Range("E33").Value = 0.9603356216628528
$ws.Range("F33").Value = 0.957564378312695
$ws.Range("C34").Value = 0.9780701754385965
$ws.Range("D34").Value = 0.9707901173606633
$ws.Range("E34").Value = 0.9780701754385965
$ws.Range("F34").Value = 0.9693926683477269
$ws.Range("C35").Value = 0.948512585812357
$ws.Range("D35").Value = 0.9343987420019269
$ws.Range("E35").Value = 0.948512585812357
$ws.Range("F35").Value = 0.9255903200972527
$ws.Range("C36").Value = 0.8771929824561403
$ws.Range("D36").Value = 0.8700135557593731
$ws.Range("E36").Value = 0.8771929824561403
$ws.Range("F36").Value = 0.8663009766207301
